$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "华工科技"
$ws.Range("B2").Value = "华工科技"
$ws.Range("C2").Value = "华工科技"
$ws.Range("A3").Value = "天奇股份"
$ws.Range("B3").Value = "电广传媒"
$ws.Range("C3").Value = "风语筑"
$ws.Range("B4").Value = "天奇股份"
$ws.Range("A5").Value = "风语筑"
$ws.Range("B5").Value = "华胜天成"
$ws.Range("C5").Value = "博纳影业"
$ws.Range("B6").Value = "风语筑"
$ws.Range("C6").Value = "天奇股份"
$ws.Range("A7").Value = "光线传媒"
$ws.Range("C7").Value = "嘉美包装"
$ws.Range("C8").Value = "利欧股份"
$ws.Range("A9").Value = "深科技"
$ws.Range("B9").Value = "东方财富"
$ws.Range("C9").Value = "光线传媒"
$ws.Range("A10").Value = "汉缆股份"
$ws.Range("B10").Value = "光线传媒"
$ws.Range("C10").Value = "汉缆股份"
$ws.Range("A11").Value = "嘉美包装"
$ws.Range("B11").Value = "深科技"
$ws.Range("C11").Value = "三花智控"
$ws.Range("A12").Value = "掌阅科技"
$ws.Range("B12").Value = "利欧股份"
$ws.Range("C12").Value = "深科技"
$ws.Range("A13").Value = "电广传媒"
$ws.Range("B13").Value = "捷成股份"
$ws.Range("C13").Value = "掌阅科技"
$ws.Range("A14").Value = "博纳影业"
$ws.Range("B14").Value = "百达精工"
$ws.Range("C14").Value = "浙江世宝"
$ws.Range("A15").Value = "百达精工"
$ws.Range("B15").Value = "汉缆股份"
$ws.Range("C15").Value = "万向钱潮"
$ws.Range("A16").Value = "万向钱潮"
$ws.Range("B16").Value = "掌阅科技"
$ws.Range("C16").Value = "协鑫集成"
$ws.Range("A17").Value = "贵州茅台"
$ws.Range("B17").Value = "万向钱潮"
$ws.Range("C17").Value = "大位科技"
$ws.Range("A18").Value = "东方财富"
$ws.Range("B18").Value = "利亚德"
$ws.Range("A19").Value = "天龙集团"
$ws.Range("B19").Value = "嘉美包装"
$ws.Range("C19").Value = "紫金矿业"
$ws.Range("A20").Value = "贵州轮胎"
$ws.Range("B20").Value = "蓝色光标"
$ws.Range("C20").Value = "百达精工"
$ws.Range("A21").Value = "拉卡拉"
$ws.Range("B21").Value = "网宿科技"
$ws.Range("C21").Value = "蓝色光标"